$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = "llama-3-8b-bnb-4bit-synthetic_text_to_sql-lora-3epochs-Q5_K_M:latest"
$ws.Cells.Item($row, 2).Value = "llama3:70b"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 200
$ws.Cells.Item($row, 5).Value = 2615.66

# F..AH (columns 6 to 34) are blank placeholder cells in this log - written as
# empty text values (same convention as the rows above it).
for ($c = 6; $c -le 34; $c++) {
    $ws.Cells.Item($row, $c).Value = "'"
    $ws.Cells.Item($row, $c).Style = "Normal"
}

$ws.Cells.Item($row, 35).Value = 2056.28  # AI

# AJ..AP (columns 36 to 42) blank placeholder cells
for ($c = 36; $c -le 42; $c++) {
    $ws.Cells.Item($row, $c).Value = "'"
    $ws.Cells.Item($row, $c).Style = "Normal"
}

$ws.Cells.Item($row, 43).Value = 559.38  # AQ
$ws.Cells.Item($row, 44).Value = 70      # AR
$ws.Cells.Item($row, 45).Value = "logs\llama_3_8b_bnb_4bit_synthetic_text_to_sql_lora_3epochs_Q5_K_M_latest_llama3_70b_1_200_test_bootstrap_match_1.txt"  # AS
$ws.Cells.Item($row, 46).Value = 559.38  # AT
$ws.Cells.Item($row, 47).Value = 35      # AU
$ws.Cells.Item($row, 48).Value = "logs\llama_3_8b_bnb_4bit_synthetic_text_to_sql_lora_3epochs_Q5_K_M_latest_llama3_70b_1_200_test_bootstrap_correct_1.txt"  # AV
$ws.Cells.Item($row, 49).Value = 58.33333333333334  # AW
$ws.Cells.Item($row, 50).Value = 4   # AX
$ws.Cells.Item($row, 51).Value = 8   # AY
